# Update the cryptos worksheet with the latest scraped prices/volumes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new D value, new E value). $null means "leave D unchanged".
$updates = @{
    2  = @("64.269.15", "  -0.05%  ")
    3  = @("3.491.15", "  -0.11%  ")
    4  = @("1.00", "  +0.06%  ")
    5  = @("587.82", "  +0.11%  ")
    6  = @("133.86", "  -0.34%  ")
    7  = @($null, "  -0.01%  ")
    8  = @("0.486", "  +0.11%  ")
    9  = @($null, "  +5.83%  ")
    10 = @("0.124", "  -0.09%  ")
    11 = @("0.388", "  +2.94%  ")
    12 = @("4.089.40", "  -0.03%  ")
    13 = @($null, "  +0.76%  ")
    14 = @($null, "  -0.18%  ")
    15 = @("3.491.69", "  -0.13%  ")
    16 = @("64.280.92", "  -0.10%  ")
    17 = @("25.32", "  +0.59%  ")
    18 = @($null, "  +1.68%  ")
    19 = @("5.76", "  +0.33%  ")
    20 = @($null, "  -0.57%  ")
    21 = @("387.60", "  -0.23%  ")
    22 = @("0.579", "  +2.53%  ")
    23 = @("3.629.64", "  -0.16%  ")
    24 = @("74.14", "  -0.39%  ")
    25 = @($null, "  +0.07%  ")
    26 = @($null, "  -0.53%  ")
    27 = @($null, "  +1.79%  ")
    28 = @($null, "  +0.02%  ")
    29 = @("7.36", "  -0.52%  ")
    30 = @($null, "  -2.40%  ")
    31 = @("2.25", "  +0.78%  ")
    32 = @("8.15", "  -1.59%  ")
    33 = @("0.154", "  +3.92%  ")
    34 = @("3.519.31", "  +0.06%  ")
    35 = @($null, "  +0.02%  ")
    36 = @("23.30", "  -0.75%  ")
    37 = @("5.32", "  +1.34%  ")
    38 = @("6.91", "  +0.50%  ")
    39 = @($null, "  +0.38%  ")
    40 = @($null, "  +2.21%  ")
    41 = @("0.0785", "  +0.45%  ")
    42 = @("0.807", "  +0.29%  ")
    43 = @($null, "  +0.09%  ")
    44 = @("4.42", "  +0.51%  ")
    45 = @("24.62", "  -3.43%  ")
    46 = @($null, "  -0.45%  ")
    47 = @("1.64", "  -0.78%  ")
    48 = @("6.80", "  +0.73%  ")
    51 = @("0.0258", "  -0.78%  ")
}

# Helper: these price/volume columns are stored as plain text in the
# workbook (t="inlineStr"). Excel will silently reinterpret a plain
# digits-and-single-dot string (e.g. "1.00", "587.82") as a number, so such
# values are written with a leading apostrophe to force text, matching the
# original text storage. Strings that cannot be parsed as a plain number
# (e.g. "64.269.15", URLs, coin names, the "  +0.06%  " volume strings) are
# already safe to set directly.
function Set-TextValue($cell, $text) {
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($null -ne $dVal) {
        Set-TextValue $ws.Cells.Item($row, 4) $dVal
    }
    Set-TextValue $ws.Cells.Item($row, 5) $eVal
}

# Rows 49 and 50 swap coin identity (Sui <-> Maker) along with new price/volume data.
Set-TextValue $ws.Cells.Item(49, 2) "Maker"
Set-TextValue $ws.Cells.Item(49, 3) "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Cells.Item(49, 4) "2.402.51"
Set-TextValue $ws.Cells.Item(49, 5) "  -3.03%  "

Set-TextValue $ws.Cells.Item(50, 2) "SuiNetwork"
Set-TextValue $ws.Cells.Item(50, 3) "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Cells.Item(50, 4) "0.918"
Set-TextValue $ws.Cells.Item(50, 5) "  +2.34%  "
